# Add "Publishing to Azure" interview-prep content:
# three new worksheets appended after the existing four:
#   "Собеседование" (interview questions), "OOP", "Sheet3" (blank, becomes active tab).

$wb = $excel.ActiveWorkbook
$nl = [char]10

$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)

# ---------------------------------------------------------------------------
# Sheet: Собеседование
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "Собеседование"

# --- Row 1 ---
$ws5.Range("A1").Value = 'Что такое кластеризованный и некластеризованный индекс?' + $nl + 'Когда какое надо использовать?'
$ws5.Range("B1").Value = 'http://sql-ex.ru/blogs/optimization/indexes_usage.html'
$ws5.Hyperlinks.Add($ws5.Range("B1"), 'http://sql-ex.ru/blogs/optimization/indexes_usage.html') | Out-Null

# --- Row 2 ---
$ws5.Range("A2").Value = 'Что такое Join? Чем он отличается от Left Join, Right Join? Inner Join? Outer Join?'
$ws5.Range("B2").Value = 'https://shra.ru/2017/09/sql-join-v-primerakh-s-opisaniem/'
$ws5.Hyperlinks.Add($ws5.Range("B2"), 'https://shra.ru/2017/09/sql-join-v-primerakh-s-opisaniem/') | Out-Null

# --- Row 3 ---
$a3 = @(
  'Есть три таблицы:'
  'CUSTOMERS (ID, NAME, MANAGER_ID);'
  'MANAGERS (ID, NAME);'
  'ORDERS (ID, DATE, AMOUNT, CUSTOMER_ID).'
  'Написать запрос, который выведет имена Customers и их SalesManagers, которые сделали покупок на общую сумму больше 10000 с 01.01.2013.'
) -join $nl
$ws5.Range("A3").Value = $a3

# --- Row 4 ---
$a4 = @(
  'Делаем электронный справочник по книгам. Ищем:'
  'А) В каком магазине купить данную книгу.'
  'Б) В каких магазинах купить книги этого автора (авторов).'
  'В) Кто автор книги.'
  'Г) Какие книги написал автор.'
  'Нарисовать БД. Написать запрос Б. (Не забыть учесть, что у одной книжки — может быть несколько авторов)'
) -join $nl
$ws5.Range("A4").Value = $a4

# --- Row 5 ---
$a5 = @(
  'Что такое агрегирующие функции?'
  'Операторы Group By, Having?'
  'Приведите примеры их использования.'
) -join $nl
$ws5.Range("A5").Value = $a5
$ws5.Range("B5").Value = 'http://dspace.ut.ee/bitstream/handle/10062/10137/_4.html;jsessionid=FC85CA23C239FC2944DE1E3780127E53'
$ws5.Hyperlinks.Add($ws5.Range("B5"), 'http://dspace.ut.ee/bitstream/handle/10062/10137/_4.html;jsessionid=FC85CA23C239FC2944DE1E3780127E53') | Out-Null
$ws5.Range("C5").Value = 'https://metanit.com/sql/sqlserver/5.2.php'
$ws5.Hyperlinks.Add($ws5.Range("C5"), 'https://metanit.com/sql/sqlserver/5.2.php') | Out-Null

# --- Row 6 ---
$a6 = @(
  'Table «PC» (id, cpu(MHz), memory(Mb), hdd(Gb))'
  '1) Тактовые частоты CPU тех компьютеров, у которых объем памяти 3000 Мб. Вывод: id, cpu, memory.'
  '2) Минимальный объём жесткого диска, установленного в компьютере на складе. Вывод: hdd.'
  '3) Количество компьютеров с минимальным объемом жесткого диска, доступного на складе. Вывод: count, hdd.'
) -join $nl
$ws5.Range("A6").Value = $a6

# --- Row 7 ---
$a7 = @(
  'Дана следующая структура базы данных в MS SQL:'
  'Departments (Id, Name), Employees(Id, DepartmentId, Name, Salary).'
  'Необходимо:'
  '• Написать запрос получения имени одного сотрудника, имеющего максимальную зарплату в компании, и название его отдела.'
  '• Получить список отделов, средняя зарплата в которых больше 1000$.'
) -join $nl
$ws5.Range("A7").Value = $a7

# --- Row 8 ---
$ws5.Range("A8").Value = 'Ado Net – что за технология? и как и когда она используется?'
$ws5.Range("B8").Value = 'https://metanit.com/sharp/adonet/1.1.php'
$ws5.Hyperlinks.Add($ws5.Range("B8"), 'https://metanit.com/sharp/adonet/1.1.php') | Out-Null

# --- Row 9 ---
$a9 = @(
  'Что такое Entity Framework?'
  'Какие подходы проектирования БД знаете?'
  'Расскажите про Code First.'
) -join $nl
$ws5.Range("A9").Value = $a9
$ws5.Range("B9").Value = 'https://www.internet-technologies.ru/articles/chto-takoe-entity-framework.html'
$ws5.Hyperlinks.Add($ws5.Range("B9"), 'https://www.internet-technologies.ru/articles/chto-takoe-entity-framework.html') | Out-Null
$ws5.Range("D9").Value = 'https://metanit.com/sharp/entityframework/1.2.php'
$ws5.Hyperlinks.Add($ws5.Range("D9"), 'https://metanit.com/sharp/entityframework/1.2.php') | Out-Null

# --- Fonts / styling ---
# NOTE: this engine only honors formatting applied to the *first* area of a
# multi-area (comma) Range union, so apply styles in per-cell loops instead.

# "Question" cells (big Georgia 18 headings), left/center, wrapped: A1,A3,A4,A6,A7
foreach ($ref in @("A1", "A3", "A4", "A6", "A7")) {
  $r = $ws5.Range($ref)
  $r.Font.Name = "Georgia"
  $r.Font.Size = 18
  $r.HorizontalAlignment = -4131
  $r.VerticalAlignment = -4108
  $r.WrapText = $true
}

# A2 Georgia 18 left/center, no wrap
$ws5.Range("A2").Font.Name = "Georgia"
$ws5.Range("A2").Font.Size = 18
$ws5.Range("A2").HorizontalAlignment = -4131
$ws5.Range("A2").VerticalAlignment = -4108

# Link cells styled like the built-in Hyperlink style, Times New Roman 15, centered: B1,B2,B5,C1(blank),C2(blank)
foreach ($ref in @("B1", "B2", "B5", "C1", "C2")) {
  $r = $ws5.Range($ref)
  $r.Font.Name = "Times New Roman"
  $r.Font.Size = 15
  $r.HorizontalAlignment = -4108
  $r.VerticalAlignment = -4108
}

# C5, B8 keep the plain Hyperlink style but centered
foreach ($ref in @("C5", "B8")) {
  $r = $ws5.Range($ref)
  $r.HorizontalAlignment = -4108
  $r.VerticalAlignment = -4108
}

# A5 Arial 16 left/center wrap
$ws5.Range("A5").Font.Name = "Arial"
$ws5.Range("A5").Font.Size = 16
$ws5.Range("A5").HorizontalAlignment = -4131
$ws5.Range("A5").VerticalAlignment = -4108
$ws5.Range("A5").WrapText = $true

# A8 Arial 16 (dark gray), no alignment override
$ws5.Range("A8").Font.Name = "Arial"
$ws5.Range("A8").Font.Size = 16
$ws5.Range("A8").Font.Color = 2236962

# A9 Arial 16 (dark gray) + wrap
$ws5.Range("A9").Font.Name = "Arial"
$ws5.Range("A9").Font.Size = 16
$ws5.Range("A9").Font.Color = 2236962
$ws5.Range("A9").WrapText = $true

# --- Row heights ---
$ws5.Rows.Item(1).RowHeight = 48
$ws5.Rows.Item(2).RowHeight = 23
$ws5.Rows.Item(3).RowHeight = 144
$ws5.Rows.Item(4).RowHeight = 144
$ws5.Rows.Item(5).RowHeight = 63
$ws5.Rows.Item(6).RowHeight = 96
$ws5.Rows.Item(7).RowHeight = 168
$ws5.Rows.Item(8).RowHeight = 20
$ws5.Rows.Item(9).RowHeight = 63

# --- Column widths ---
$ws5.Columns.Item(1).ColumnWidth = 163.5
$ws5.Columns.Item(2).ColumnWidth = 122.33
$ws5.Columns.Item(3).ColumnWidth = 122.33
$ws5.Columns.Item(4).ColumnWidth = 42.67

# --- View state ---
$ws5.Range("D9").Select() | Out-Null

Write-Host "Собеседование sheet populated"
